$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster (A2) changes from ECs to MuSCs
$ws.Range("A2").Value = "MuSCs"

# Target cluster (D2) stays ECs (no visible change, but forces re-resolution)
$ws.Range("D2").Value = "ECs"

# Updated TPM-derived numeric values for row 2
$ws.Range("G2").Value = 0.2195956666666667
$ws.Range("H2").Value = 0.658787
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.009232000000000001
$ws.Range("N2").Value = 0.027696
$ws.Range("Q2").Value = 0.002027307194666667
$ws.Range("R2").Value = 0.018245764752
